$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The cell A3 held the "${String}" placeholder (a shared string). Remove
# its content so the cell becomes blank again; this also drops the now
# unused entry from the shared strings table on save.
$ws.Range("A3").ClearContents()
$ws.Range("A3").Style = "Normal"

# Update the active selection to G10, matching the worksheet view saved
# in the target workbook.
[void]$ws.Range("G10").Select()
